$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D slightly to fit the new, longer wrapped comment text
$ws.Columns.Item(4).ColumnWidth = 27.8

# --- Row 5 ("SNO" 3): Databinding ngFor,styleManagement / completed ---
$ws.Cells.Item(4, 2).Copy()
$ws.Cells.Item(5, 2).PasteSpecial(-4122)
$ws.Cells.Item(3, 4).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4122)

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 44261
$ws.Cells.Item(5, 3).Value = "8hr"
$ws.Cells.Item(5, 4).Value = "Databinding ngFor,styleManagement"
$ws.Cells.Item(5, 5).Value = "completed"
$ws.Rows.Item(5).RowHeight = 30

# --- Row 6 ("SNO" 4): Modules & Services / not completed ---
$ws.Cells.Item(4, 2).Copy()
$ws.Cells.Item(6, 2).PasteSpecial(-4122)

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 44262
$ws.Cells.Item(6, 3).Value = "8hr"
$ws.Cells.Item(6, 4).Value = "Modules & Services"
$ws.Cells.Item(6, 5).Value = "not completed"

$ws.Application.CutCopyMode = $false

$ws.Range("D11").Select()
